$wb = $excel.ActiveWorkbook

$csa = $wb.Worksheets.Item("csa")
$eee = $wb.Worksheets.Item("eee")

# --- Add the new "ece" worksheet right after "eee" ---
$ece = $wb.Worksheets.Add($null, $eee)
$ece.Name = "ece"

# --- Populate ece with its timetable data ---
$ece.Range("A1").Value = "date"
$ece.Range("B1").Value = "subject"
$ece.Range("A2").Value = 45092
$ece.Range("B2").Value = "electricity"
$ece.Range("A3").Value = 45097
$ece.Range("B3").Value = "kseb"

# Copy the date number format from an existing sheet so the new date cells
# reuse the workbook's existing style (instead of creating a new one).
$csa.Range("A2:A3").Copy()
$ece.Range("A2:A3").PasteSpecial(-4122)   # xlPasteFormats

# --- Update selection on "eee" (no longer the active tab) ---
$eee.Select()
$eee.Range("A1:B1").Select()

# --- Finally activate "ece" as the selected/visible tab ---
$ece.Select()
$ece.Range("C8").Select()
